$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# workbook.xml: calcPr gains refMode="R1C1" -> switch Excel's reference style to R1C1
$excel.ReferenceStyle = -4150  # xlR1C1

# sheet1.xml: widen column C from 38.42578125 to 51 (character width units)
$ws.Columns.Item(3).ColumnWidth = 50.166666666666664

# sheet1.xml: swap the text shown in the merged C1:G1 / C2:G2 banner cells
$tmp = $ws.Range("C1").Value2
$ws.Range("C1").Value2 = $ws.Range("C2").Value2
$ws.Range("C2").Value2 = $tmp

# sheet1.xml: move the active selection from K5 to J4
$ws.Range("J4").Select()
